$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'272.70"
$ws.Range("E2").Value = "'4.34%"
$ws.Range("G2").Value = "'23"

$ws.Range("E3").Value = "'-1.44%"
$ws.Range("G3").Value = "'23"

$ws.Range("D4").Value = "'4.728"
$ws.Range("E4").Value = "'0.30%"
$ws.Range("G4").Value = "'23"

$ws.Range("D5").Value = "'0.06131"
$ws.Range("E5").Value = "'-1.19%"
$ws.Range("G5").Value = "'23"

$ws.Range("D6").Value = "'6.744"
$ws.Range("E6").Value = "'0.15%"
$ws.Range("G6").Value = "'23"

$ws.Range("D7").Value = "'0.8552"
$ws.Range("E7").Value = "'0.40%"
$ws.Range("G7").Value = "'23"

$ws.Range("D8").Value = "'0.9046"
$ws.Range("E8").Value = "'-0.35%"
$ws.Range("G8").Value = "'23"

$ws.Range("D9").Value = "'0.1435"
$ws.Range("E9").Value = "'2.18%"
$ws.Range("G9").Value = "'23"

$ws.Range("D10").Value = "'0.05044"
$ws.Range("E10").Value = "'2.44%"
$ws.Range("G10").Value = "'23"

$ws.Range("D11").Value = "'0.07156"
$ws.Range("G11").Value = "'23"

$ws.Range("D12").Value = "'0.03171"
$ws.Range("E12").Value = "'-0.12%"
$ws.Range("G12").Value = "'23"

$ws.Range("D13").Value = "'0.09036"
$ws.Range("E13").Value = "'-0.23%"
$ws.Range("G13").Value = "'23"

$ws.Range("D14").Value = "'0.001545"
$ws.Range("E14").Value = "'0.06%"
$ws.Range("G14").Value = "'23"

$ws.Range("D15").Value = "'0.0006082"
$ws.Range("E15").Value = "'-1.28%"
$ws.Range("G15").Value = "'23"

$ws.Range("D16").Value = "'0.005942"
$ws.Range("E16").Value = "'-2.45%"
$ws.Range("G16").Value = "'23"

$ws.Range("D17").Value = "'3.461"
$ws.Range("E17").Value = "'-0.19%"
$ws.Range("G17").Value = "'23"

$ws.Range("D18").Value = "'3.186"
$ws.Range("E18").Value = "'0.49%"
$ws.Range("G18").Value = "'23"

$ws.Range("D19").Value = "'2.263"
$ws.Range("E19").Value = "'3.95%"
$ws.Range("G19").Value = "'23"

$ws.Range("E20").Value = "'-0.70%"
$ws.Range("G20").Value = "'23"

$ws.Range("E21").Value = "'0.06%"
$ws.Range("G21").Value = "'23"

$ws.Range("D22").Value = "'3.833"
$ws.Range("E22").Value = "'-6.72%"
$ws.Range("G22").Value = "'23"

$ws.Range("D23").Value = "'0.04241"
$ws.Range("E23").Value = "'-0.10%"
$ws.Range("G23").Value = "'23"

$ws.Range("D24").Value = "'0.001179"
$ws.Range("E24").Value = "'-3.20%"
$ws.Range("G24").Value = "'23"

$ws.Range("D25").Value = "'0.004154"
$ws.Range("E25").Value = "'0.89%"
$ws.Range("G25").Value = "'23"

$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'-0.13%"
$ws.Range("G26").Value = "'23"

$ws.Range("E27").Value = "'3.84%"
$ws.Range("G27").Value = "'23"

$ws.Range("G28").Value = "'23"

$ws.Range("G29").Value = "'23"

$ws.Range("G30").Value = "'23"

$ws.Range("G31").Value = "'23"

$ws.Range("G32").Value = "'23"

$ws.Range("G33").Value = "'23"

$ws.Range("G34").Value = "'23"

$ws.Range("G35").Value = "'23"

$ws.Range("G36").Value = "'23"

$ws.Range("G37").Value = "'23"

$ws.Range("G38").Value = "'23"

$ws.Range("G39").Value = "'23"

$ws.Range("D40").Value = "'0.03973"
$ws.Range("E40").Value = "'1.73%"
$ws.Range("G40").Value = "'23"

$ws.Range("E41").Value = "'0.58%"
$ws.Range("G41").Value = "'23"

$ws.Range("D42").Value = "'0.004198"
$ws.Range("E42").Value = "'1.54%"
$ws.Range("G42").Value = "'23"

$ws.Range("D43").Value = "'0.002088"
$ws.Range("E43").Value = "'-4.39%"
$ws.Range("G43").Value = "'23"

$ws.Range("D44").Value = "'0.01178"
$ws.Range("E44").Value = "'-12.18%"
$ws.Range("G44").Value = "'23"

$ws.Range("D45").Value = "'0.00005136"
$ws.Range("E45").Value = "'-0.77%"
$ws.Range("G45").Value = "'23"

$ws.Range("E46").Value = "'0.06%"
$ws.Range("G46").Value = "'23"

$ws.Range("D47").Value = "'0.8976"
$ws.Range("E47").Value = "'1,452.36%"
$ws.Range("G47").Value = "'23"

$ws.Range("G48").Value = "'23"

$ws.Range("E49").Value = "'0.06%"
$ws.Range("G49").Value = "'23"

$ws.Range("E50").Value = "'0.06%"
$ws.Range("G50").Value = "'23"

$ws.Range("G51").Value = "'23"
